$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (shifts existing row 2 data down to row 3)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Force column A to be treated as text so the date-like string "2026-01-10"
# is not auto-converted into a date serial number.
$ws.Range("A2").NumberFormat = "@"

# Populate the new row 2 with the latest news item
$ws.Range("A2").Value = "2026-01-10"
$ws.Range("B2").Value = "Estudiantes respiran con el aumento del pasaje de TransMilenio: no tendrán que gastarse lo del almuerzo"
$ws.Range("C2").Value = "Alerta Bogotá"
$ws.Range("D2").Value = "Bogotá"
$ws.Range("E2").Value = "https://www.alertabogota.com/noticias/local/estudiantes-respiran-con-el-aumento-del-pasaje-de-transmilenio-no-tendran-que-gastarse-lo-del-almuerzo"
$ws.Range("F2").Value = "El incremento en el pasaje no será una barrera para que los estudiantes continúen asistiendo a clases."

# Drop the temporary number-format override so the cell carries no explicit
# style, matching the rest of the data rows.
$ws.Range("A2").ClearFormats()

Write-Output "done"
